# edit.ps1
# Applies the "Updated cryptos list" data refresh to Sheet1 of the cryptos
# workbook: refreshed Price (column D) and Volume(1h) (column E) values for
# most rows, plus two row-pairs whose Coin/Link/Price/Volume were swapped
# (rows 30<->31 and 46<->47) because the ranking order changed.
#
# Cells in column D can look like plain numbers (e.g. "1.00", "0.480").
# Excel's COM layer would silently reinterpret such a string as a number
# and drop the formatting (e.g. "1.00" -> 1, "0.480" -> 0.48), so for those
# cells we force the cell's number format to Text ("@") before assigning
# the value, guaranteeing the literal string is preserved exactly as in
# the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet



$ws.Cells.Item(2,4).Value = '62.180.62'
$ws.Cells.Item(2,5).Value = '  -4.19%  '

$ws.Cells.Item(3,4).Value = '3.026.22'
$ws.Cells.Item(3,5).Value = '  -3.56%  '

$ws.Cells.Item(4,5).Value = '  +0.29%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '536.62'
$ws.Cells.Item(5,5).Value = '  -4.95%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '131.60'
$ws.Cells.Item(6,5).Value = '  -11.22%  '

$ws.Cells.Item(7,5).Value = '  +0.19%  '

$ws.Cells.Item(8,4).Value = '3.019.10'
$ws.Cells.Item(8,5).Value = '  -3.48%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.480'
$ws.Cells.Item(9,5).Value = '  -3.55%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '6.32'
$ws.Cells.Item(10,5).Value = '  -8.44%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.152'
$ws.Cells.Item(11,5).Value = '  -3.85%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.450'
$ws.Cells.Item(12,5).Value = '  -2.88%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '33.98'
$ws.Cells.Item(13,5).Value = '  -5.23%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.0000210'
$ws.Cells.Item(14,5).Value = '  -5.41%  '

$ws.Cells.Item(15,4).Value = '3.524.46'
$ws.Cells.Item(15,5).Value = '  -3.10%  '

$ws.Cells.Item(16,4).Value = '62.442.80'
$ws.Cells.Item(16,5).Value = '  -3.77%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.110'
$ws.Cells.Item(17,5).Value = '  -2.19%  '

$ws.Cells.Item(18,4).Value = '3.036.54'
$ws.Cells.Item(18,5).Value = '  -3.18%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '6.48'
$ws.Cells.Item(19,5).Value = '  -3.76%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '474.19'
$ws.Cells.Item(20,5).Value = '  -9.74%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '13.11'
$ws.Cells.Item(21,5).Value = '  -5.24%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '0.683'
$ws.Cells.Item(22,5).Value = '  -2.53%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '6.92'
$ws.Cells.Item(23,5).Value = '  -7.10%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '76.20'
$ws.Cells.Item(24,5).Value = '  -2.91%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '11.92'
$ws.Cells.Item(25,5).Value = '  -6.09%  '

$ws.Cells.Item(26,5).Value = '  +0.00%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '2.64'
$ws.Cells.Item(27,5).Value = '  -5.24%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '8.09'
$ws.Cells.Item(28,5).Value = '  -6.21%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '1.00'
$ws.Cells.Item(29,5).Value = '  -0.01%  '

$ws.Cells.Item(30,2).Value = 'EthereumClassic'
$ws.Cells.Item(30,3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '25.69'
$ws.Cells.Item(30,5).Value = '  -1.60%  '

$ws.Cells.Item(31,2).Value = 'ImmutableX'
$ws.Cells.Item(31,3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '1.88'
$ws.Cells.Item(31,5).Value = '  -11.21%  '

$ws.Cells.Item(32,5).Value = '  -3.88%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '59.27'
$ws.Cells.Item(33,5).Value = '  +12.28%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '2.42'
$ws.Cells.Item(34,5).Value = '  -8.15%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '503.62'
$ws.Cells.Item(35,5).Value = '  -10.48%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '5.77'
$ws.Cells.Item(36,5).Value = '  -4.32%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '5.01'
$ws.Cells.Item(37,5).Value = '  -6.72%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.0390'
$ws.Cells.Item(38,5).Value = '  -11.05%  '

$ws.Cells.Item(39,4).Value = '3.029.47'
$ws.Cells.Item(39,5).Value = '  -0.94%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.0772'
$ws.Cells.Item(40,5).Value = '  -5.05%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.115'
$ws.Cells.Item(41,5).Value = '  -4.60%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '7.89'
$ws.Cells.Item(42,5).Value = '  -4.23%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '2.53'
$ws.Cells.Item(43,5).Value = '  -10.61%  '

$ws.Cells.Item(44,5).Value = '  +0.03%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.246'
$ws.Cells.Item(45,5).Value = '  -3.86%  '

$ws.Cells.Item(46,2).Value = 'Monero'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '119.45'
$ws.Cells.Item(46,5).Value = '  -0.25%  '

$ws.Cells.Item(47,2).Value = 'Fetch.AI'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '1.97'
$ws.Cells.Item(47,5).Value = '  -8.74%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '23.64'
$ws.Cells.Item(48,5).Value = '  -5.32%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '0.105'
$ws.Cells.Item(49,5).Value = '  -3.04%  '

$ws.Cells.Item(50,5).Value = '  +60.15%  '

$ws.Cells.Item(51,4).Value = '0.0₃0481'
$ws.Cells.Item(51,5).Value = '  -7.96%  '
